$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Page Size" notes cell (H6): Min=256 -> Min=512
$ws.Range("H6").Value = "Can only be changed when the database is empty (before the first CREATE TABLE statement. Min=512, Max=65536 "

# Update the "Secure Delete" row (row 12): Is verifiable? and Persistent? columns
$ws.Range("E12").Value = "Yes"
$ws.Range("G12").Value = "No"

# Update the selected cell/active cell to D6
$ws.Range("D6").Select()
